# Macedonian_data.xlsx edit script
# Renames the bare "NP" argument label (the unmarked / default-case argument,
# which is actually the grammatical subject of the Macedonian clause) to "SBJ"
# throughout the valency table, and updates the corresponding valency_pattern
# strings (column L) to match. Also normalizes the previously-blank Y column
# for plain transitive ("TR") rows to "ACC", and fills in "*" for the single
# irregular row. Finally resets the view's selection to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (language_no) starting at row 2.
$lastRow = $ws.Cells(1, 1).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $xVal = $ws.Cells.Item($r, 9).Value()   # column I -> X
    $yVal = $ws.Cells.Item($r, 10).Value()  # column J -> Y
    $locus = $ws.Cells.Item($r, 11).Value() # column K -> locus
    $pattern = $ws.Cells.Item($r, 12).Value() # column L -> valency_pattern

    $newX = $xVal
    $newY = $yVal
    $newPattern = $pattern

    if ($locus -eq "X") {
        # X is the marked argument; Y is the bare/default one -> becomes SBJ
        if ($yVal -eq "NP") {
            $newY = "SBJ"
        }
        if ($pattern) {
            $parts = $pattern.Split("_")
            if ($parts.Count -eq 2 -and $parts[1] -eq "NP") {
                $newPattern = $parts[0] + "_SBJ"
            }
        }
    }
    elseif ($locus -eq "Y") {
        # Y is the marked argument; X is the bare/default one -> becomes SBJ
        if ($xVal -eq "NP") {
            $newX = "SBJ"
        }
        if ($pattern) {
            $parts = $pattern.Split("_")
            if ($parts.Count -ge 1 -and $parts[0] -eq "NP") {
                $parts[0] = "SBJ"
                $newPattern = [string]::Join("_", $parts)
            }
        }
    }
    elseif ($locus -eq "TR") {
        # Plain transitive verb: X is the subject, Y was left blank but is
        # understood to be the accusative object.
        $newX = "SBJ"
        $newY = "ACC"
        # locus/pattern (both "TR") remain unchanged.
    }
    elseif ($locus -eq "*") {
        # Irregular row: mirror the "*" into the Y column as well.
        $newY = "*"
    }

    if ($newX -ne $xVal) {
        $ws.Cells.Item($r, 9).Value = $newX
    }
    if ($newY -ne $yVal) {
        $ws.Cells.Item($r, 10).Value = $newY
    }
    if ($newPattern -ne $pattern) {
        $ws.Cells.Item($r, 12).Value = $newPattern
    }
}

# Reset the sheet view: clear any frozen/scrolled top-left cell and select A8.
$ws.Range("A8").Select()
